$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Overview")

# --- Update copyright year ---
$ws.Range("B3").Value = "Copyright @2015 - 2023"
$ws.Range("E8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("I8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "-"
$ws.Range("E12").Value = 8837
$ws.Range("F12").Value = "-"
$ws.Range("G12").Value = 7754
$ws.Range("H12").Value = 11452
$ws.Range("I12").Value = 10757
$ws.Range("E13").Value = 10714
$ws.Range("F13").Value = "-"
$ws.Range("G13").Value = 9733
$ws.Range("H13").Value = 16486
$ws.Range("I13").Value = 16617
$ws.Range("E14").Value = 11473
$ws.Range("F14").Value = "-"
$ws.Range("G14").Value = 9152
$ws.Range("H14").Value = 13306
$ws.Range("I14").Value = 14049
$ws.Range("E15").Value = 31024
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 26639
$ws.Range("H15").Value = 41244
$ws.Range("I15").Value = 41423
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = "-"
$ws.Range("E18").Value = 63
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = 178
$ws.Range("I18").Value = 91
$ws.Range("E19").Value = 295
$ws.Range("F19").Value = "-"
$ws.Range("G19").Value = 576
$ws.Range("H19").Value = 486
$ws.Range("I19").Value = 365
$ws.Range("E20").Value = 35
$ws.Range("F20").Value = "-"
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 198
$ws.Range("I20").Value = 175
$ws.Range("E21").Value = 393
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 622
$ws.Range("H21").Value = 862
$ws.Range("I21").Value = 631
$ws.Range("F26").Value = "-"
$ws.Range("E27").Value = 31417
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 27261
$ws.Range("H27").Value = 42106
$ws.Range("I27").Value = 42054
$ws.Range("E31").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("F31").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("G31").Value = "فصل اول منتهی به 1401/03"
$ws.Range("H31").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("I31").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = "-"
$ws.Range("E35").Value = 2139558
$ws.Range("F35").Value = "-"
$ws.Range("G35").Value = 2242593
$ws.Range("H35").Value = 3483135
$ws.Range("I35").Value = 3778962
$ws.Range("E36").Value = 2000284
$ws.Range("F36").Value = "-"
$ws.Range("G36").Value = 2155377
$ws.Range("H36").Value = 3858958
$ws.Range("I36").Value = 4291091
$ws.Range("E37").Value = 1206770
$ws.Range("F37").Value = "-"
$ws.Range("G37").Value = 1258059
$ws.Range("H37").Value = 1884248
$ws.Range("I37").Value = 2058189
$ws.Range("E38").Value = 5346612
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 5656029
$ws.Range("H38").Value = 9226341
$ws.Range("I38").Value = 10128242
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = "-"
$ws.Range("E41").Value = 16464
$ws.Range("F41").Value = "-"
$ws.Range("G41").Value = 9666
$ws.Range("H41").Value = 84350
$ws.Range("I41").Value = 37306
$ws.Range("E42").Value = 44525
$ws.Range("F42").Value = "-"
$ws.Range("G42").Value = 46649
$ws.Range("H42").Value = 69197
$ws.Range("I42").Value = 73912
$ws.Range("E43").Value = 7435
$ws.Range("F43").Value = "-"
$ws.Range("G43").Value = 3445
$ws.Range("H43").Value = 34152
$ws.Range("I43").Value = 33798
$ws.Range("E44").Value = 68424
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 59760
$ws.Range("H44").Value = 187699
$ws.Range("I44").Value = 145016
$ws.Range("G46").Value = 14639
$ws.Range("H46").Value = 25600
$ws.Range("I46").Value = 11917
$ws.Range("G47").Value = 14639
$ws.Range("H47").Value = 25600
$ws.Range("I47").Value = 11917
$ws.Range("F49").Value = "-"
$ws.Range("E50").Value = 5415036
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 5730428
$ws.Range("H50").Value = 9439640
$ws.Range("I50").Value = 10285175
$ws.Range("E54").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("F54").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("G54").Value = "فصل اول منتهی به 1401/03"
$ws.Range("H54").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("I54").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("E58").Value = 237043873
$ws.Range("F58").Value = 243325015
$ws.Range("G58").Value = 289217565
$ws.Range("H58").Value = 304150803
$ws.Range("I58").Value = 351302594
$ws.Range("E59").Value = 186698152
$ws.Range("F59").Value = 190538489
$ws.Range("G59").Value = 221450426
$ws.Range("H59").Value = 234074851
$ws.Range("I59").Value = 258235000
$ws.Range("E60").Value = 105192643
$ws.Range("F60").Value = 113863096
$ws.Range("G60").Value = 137462740
$ws.Range("H60").Value = 141608898
$ws.Range("I60").Value = 146500747
$ws.Range("E63").Value = 265548387
$ws.Range("F63").Value = 322276316
$ws.Range("G63").Value = 371769231
$ws.Range("H63").Value = 473876404
$ws.Range("I63").Value = 409956044
$ws.Range("E64").Value = 150932203
$ws.Range("F64").Value = 94268546
$ws.Range("G64").Value = 80987847
$ws.Range("H64").Value = 142380658
$ws.Range("I64").Value = 202498630
$ws.Range("E65").Value = 212428571
$ws.Range("F65").Value = 175292683
$ws.Range("G65").Value = 172250000
$ws.Range("H65").Value = 172484848
$ws.Range("I65").Value = 193131429
$ws.Range("E71").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("F71").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("G71").Value = "فصل اول منتهی به 1401/03"
$ws.Range("H71").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("I71").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = "-"
$ws.Range("E75").Value = -1517322
$ws.Range("F75").Value = -2046685
$ws.Range("G75").Value = -1504719
$ws.Range("H75").Value = -2227176
$ws.Range("I75").Value = -2311877
$ws.Range("E76").Value = -1637348
$ws.Range("F76").Value = -2153379
$ws.Range("G76").Value = -1710430
$ws.Range("H76").Value = -2894705
$ws.Range("I76").Value = -3020615
$ws.Range("E77").Value = -1037358
$ws.Range("F77").Value = -1323228
$ws.Range("G77").Value = -1002838
$ws.Range("H77").Value = -1443000
$ws.Range("I77").Value = -1513782
$ws.Range("E78").Value = -4192028
$ws.Range("F78").Value = -5523292
$ws.Range("G78").Value = -4217987
$ws.Range("H78").Value = -6564881
$ws.Range("I78").Value = -6846274
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = "-"
$ws.Range("E81").Value = -12467
$ws.Range("F81").Value = -19560
$ws.Range("G81").Value = -6931
$ws.Range("H81").Value = -63834
$ws.Range("I81").Value = 70765
$ws.Range("E82").Value = -38163
$ws.Range("F82").Value = -63174
$ws.Range("G82").Value = -45766
$ws.Range("H82").Value = -57809
$ws.Range("I82").Value = 103575
$ws.Range("E83").Value = -5784
$ws.Range("F83").Value = -10669
$ws.Range("G83").Value = -3043
$ws.Range("H83").Value = -30672
$ws.Range("I83").Value = 33715
$ws.Range("E84").Value = -56414
$ws.Range("F84").Value = -93403
$ws.Range("G84").Value = -55740
$ws.Range("H84").Value = -152315
$ws.Range("I84").Value = 208055
$ws.Range("F86").Value = -33207
$ws.Range("G86").Value = -13000
$ws.Range("H86").Value = -14458
$ws.Range("I86").Value = -6493
$ws.Range("F87").Value = -33207
$ws.Range("G87").Value = -13000
$ws.Range("H87").Value = -14458
$ws.Range("I87").Value = -6493
$ws.Range("E90").Value = -4248442
$ws.Range("F90").Value = -5649902
$ws.Range("G90").Value = -4286727
$ws.Range("H90").Value = -6731654
$ws.Range("I90").Value = -6644712
$ws.Range("E94").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("F94").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("G94").Value = "فصل اول منتهی به 1401/03"
$ws.Range("H94").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("I94").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = "-"
$ws.Range("E98").Value = 622236
$ws.Range("F98").Value = 810194
$ws.Range("G98").Value = 737874
$ws.Range("H98").Value = 1255959
$ws.Range("I98").Value = 1467085
$ws.Range("E99").Value = 362936
$ws.Range("F99").Value = 512445
$ws.Range("G99").Value = 444947
$ws.Range("H99").Value = 964253
$ws.Range("I99").Value = 1270476
$ws.Range("E100").Value = 169412
$ws.Range("F100").Value = 335188
$ws.Range("G100").Value = 255221
$ws.Range("H100").Value = 441248
$ws.Range("I100").Value = 544407
$ws.Range("E101").Value = 1154584
$ws.Range("F101").Value = 1657827
$ws.Range("G101").Value = 1438042
$ws.Range("H101").Value = 2661460
$ws.Range("I101").Value = 3281968
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = "-"
$ws.Range("E104").Value = 3997
$ws.Range("F104").Value = 4933
$ws.Range("G104").Value = 2735
$ws.Range("H104").Value = 20516
$ws.Range("I104").Value = 108071
$ws.Range("E105").Value = 6362
$ws.Range("F105").Value = 363
$ws.Range("G105").Value = 883
$ws.Range("H105").Value = 11388
$ws.Range("I105").Value = 177487
$ws.Range("E106").Value = 1651
$ws.Range("F106").Value = 3705
$ws.Range("G106").Value = 402
$ws.Range("H106").Value = 3480
$ws.Range("I106").Value = 67513
$ws.Range("E107").Value = 12010
$ws.Range("F107").Value = 9001
$ws.Range("G107").Value = 4020
$ws.Range("H107").Value = 35384
$ws.Range("I107").Value = 353071
$ws.Range("F109").Value = 11828
$ws.Range("G109").Value = 1639
$ws.Range("H109").Value = 11142
$ws.Range("I109").Value = 5424
$ws.Range("F110").Value = 11828
$ws.Range("G110").Value = 1639
$ws.Range("H110").Value = 11142
$ws.Range("I110").Value = 5424
$ws.Range("E111").Value = 1166594
$ws.Range("F111").Value = 1678656
$ws.Range("G111").Value = 1443701
$ws.Range("H111").Value = 2707986
$ws.Range("I111").Value = 3640463